# 2018-02-26_shorter_names -> shorten the "Diabetes" node name to "Diabetes etc"
# and move the active-cell selection to B5 (the row that was just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Diabetes etc"

$ws.Range("B5").Select()
